# Actualizo código del CRM
# Add new client record (C1002 - Tatiana Avila) to the Clientes sheet,
# mirroring the previous row's sucursal/asesor/fechas/estatus pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 4

$ws.Range("A$newRow").Value = "C1002"
$ws.Range("B$newRow").Value = "Tatiana Avila"
$ws.Range("C$newRow").Value = "mundo E"
$ws.Range("D$newRow").Value = "Martha Ortiz"

# The fecha_ingreso / fecha_dispersion columns store plain text dates
# ("2025-10-08") rather than real date serials elsewhere in this sheet,
# so force text formatting before assigning, then restore the default
# style so no stray number-format sticks to the cell.
$ws.Range("E$newRow").NumberFormat = "@"
$ws.Range("E$newRow").Value = "2025-10-08"
$ws.Range("E$newRow").Style = "Normal"

$ws.Range("F$newRow").NumberFormat = "@"
$ws.Range("F$newRow").Value = "2025-10-08"
$ws.Range("F$newRow").Style = "Normal"

$ws.Range("G$newRow").Value = "DISPERSADO"
